# Reformat "Usernamepasswords" sheet: transpose the per-service columns into
# rows (Name / Username / Password), add Email + Email-alias rows, add an
# autofilter, and relink all the hyperlinks to the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- wipe the old layout (content, formatting and hyperlinks) -------------
$ws.Hyperlinks.Delete()
$ws.Cells.Clear()

# --- header row -------------------------------------------------------
# Bolding the whole row first (before any cell holds a value) is what makes
# the engine stamp row 1 itself with the bold style (customFormat="1"),
# matching how the original file represents the header row.
$ws.Rows.Item(1).Font.Bold = $true

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Username"
$ws.Range("C1").Value = "Password"

# --- data rows ----------------------------------------------------------
# Each former column (Twitter, Instagram, ...) becomes a row of
# Name / Username / Password, plus four new "Email" rows up top.
# (username, style-as-hyperlink?) pairs line up with $data below.
$data = @(
    @("Email", "info@womencoders.org", "Peoplespaceoc1", $true),
    @("Email Alias: angela.li@womencoders.org", "info@womencoders.org", "Peoplespaceoc1", $true),
    @("Email Alias: laurie.tran@womencoders.org", "info@womencoders.org", "Peoplespaceoc1", $true),
    @("Email Alias: melinda.kobayashi@womencoders.org", "info@womencoders.org", "Peoplespaceoc1", $true),
    @("Eventbrite", "Womencoders@gmail.com", "Peoplespaceoc", $true),
    @("Facebook", "womencoders", "(will add you girls as host)", $false),
    @("FTP Server", "womencoders", "peoplespace88", $false),
    @("Gmail", "womencoders@gmail.com", "Peoplespaceoc", $true),
    @("Hootsuite", "womencoders@gmail.com", "Peoplespaceoc", $true),
    @("Instagram", "womencoders", "peoplespaceoc", $false),
    @("Mailchimp", "womencoders@gmail.com", "Peoplespaceoc1", $true),
    @("Stripe", "womencoders@gmail.com", "Peoplespaceoc", $true),
    @("Twitter", "womencoders", "Peoplespaceoc", $false)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Leftover styled-but-empty cells in G2:I5 (mirrors the source workbook).
foreach ($rr in 2..5) {
    foreach ($cc in @(7,8,9)) {
        $ws.Cells.Item($rr, $cc).Style = "Hyperlink"
    }
}

# --- hyperlinks -----------------------------------------------------------
# (Hyperlinks.Add re-applies its own formatting to the target cell, so the
# explicit Style="Hyperlink" pass below runs afterwards to normalise it back
# to the workbook's existing Hyperlink cell style.)
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:womencoders@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:Womencoders@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B12"), "mailto:womencoders@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B13"), "mailto:womencoders@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:womencoders@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:info@womencoders.org") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:info@womencoders.org") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:info@womencoders.org") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:info@womencoders.org") | Out-Null

foreach ($rr in @(2,3,4,5,6,9,10,12,13)) {
    $ws.Cells.Item($rr, 2).Style = "Hyperlink"
}

# --- column widths --------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 29
$ws.Columns.Item(9).ColumnWidth = 18.33

# --- autofilter + defined name --------------------------------------------
$ws.Range("A1:C17").AutoFilter() | Out-Null
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$C`$14")
$n.Visible = $false

# --- sheet view: drop the old topLeftCell/selection, point at C16 --------
$ws.Range("C16").Select()

Write-Host "done"
